$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map: row number -> @{ D = "new price text (or $null if unchanged)"; E = "new volume% text" }
$updates = @{
    2  = @{ D = "65.135.07"; E = "  -1.08%  " }
    3  = @{ D = "2.941.15";  E = "  -2.45%  " }
    4  = @{ D = $null;       E = "  -0.04%  " }
    5  = @{ D = "567.47";    E = "  -3.10%  " }
    6  = @{ D = "158.41";    E = "  +2.33%  " }
    7  = @{ D = $null;       E = "  +0.02%  " }
    8  = @{ D = $null;       E = "  -0.14%  " }
    9  = @{ D = "2.938.03";  E = "  -2.42%  " }
    10 = @{ D = $null;       E = "  -4.44%  " }
    11 = @{ D = $null;       E = "  -3.01%  " }
    12 = @{ D = $null;       E = "  +1.15%  " }
    13 = @{ D = $null;       E = "  +0.70%  " }
    14 = @{ D = $null;       E = "  -0.18%  " }
    15 = @{ D = $null;       E = "  -0.80%  " }
    16 = @{ D = "65.156.45"; E = "  -1.00%  " }
    17 = @{ D = "3.430.61";  E = "  -2.46%  " }
    18 = @{ D = "6.93";      E = "  -0.62%  " }
    19 = @{ D = "2.938.75";  E = "  -2.65%  " }
    20 = @{ D = "14.88";     E = "  +7.94%  " }
    21 = @{ D = "444.98";    E = "  -3.38%  " }
    22 = @{ D = "0.686";     E = "  +0.19%  " }
    23 = @{ D = "7.22";      E = "  -1.99%  " }
    24 = @{ D = "82.13";     E = "  +0.11%  " }
    25 = @{ D = "2.20";      E = "  -2.19%  " }
    26 = @{ D = "12.06";     E = "  -4.30%  " }
    27 = @{ D = $null;       E = "  +0.06%  " }
    28 = @{ D = "10.01";     E = "  -7.28%  " }
    29 = @{ D = $null;       E = "  +0.74%  " }
    30 = @{ D = $null;       E = "  -1.68%  " }
    31 = @{ D = $null;       E = "  -1.65%  " }
    32 = @{ D = $null;       E = "  -2.53%  " }
    33 = @{ D = "27.10";     E = "  +0.36%  " }
    34 = @{ D = "0.110";     E = "  -1.63%  " }
    35 = @{ D = "0.999";     E = "  -0.03%  " }
    36 = @{ D = "0.972";     E = "  -2.05%  " }
    37 = @{ D = $null;       E = "  -1.27%  " }
    38 = @{ D = "49.59";     E = "  +0.27%  " }
    39 = @{ D = $null;       E = "  -1.89%  " }
    40 = @{ D = $null;       E = "  -8.28%  " }
    41 = @{ D = $null;       E = "  -1.14%  " }
    42 = @{ D = $null;       E = "  -5.93%  " }
    43 = @{ D = "0.297";     E = "  -1.86%  " }
    44 = @{ D = "8.45";      E = "  -0.25%  " }
    45 = @{ D = "385.31";    E = "  -1.00%  " }
    46 = @{ D = "0.0350";    E = "  -0.86%  " }
    47 = @{ D = "2.704.78";  E = "  -3.36%  " }
    48 = @{ D = "132.98";    E = "  -1.47%  " }
    50 = @{ D = $null;       E = "  +4.59%  " }
    51 = @{ D = $null;       E = "  -0.42%  " }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]

    if ($null -ne $vals.D) {
        # Some new price strings (e.g. "567.47") look like plain numbers to
        # Excel's auto-detection and would otherwise be stored as numeric
        # values instead of text. Force text formatting, write, then drop
        # back to the default (unstyled) cell style so no stray "s"
        # attribute / number format is left behind on save.
        $dCell = $ws.Cells.Item($row, 4)
        $dCell.NumberFormat = "@"
        $dCell.Value = $vals.D
        $dCell.Style = "Normal"
    }

    $ws.Cells.Item($row, 5).Value = $vals.E
}
